# alertas y correciones generales
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Actual")

$data = @(
    @("101546", "Leo", "Operaciones", "2024-01-29", "11:9:21", "11:09:22", "11:09:22"),
    @("101546", "Leo", "Operaciones", "2024-01-29", "11:13:37", "11:09:22", "11:13:37"),
    @("12345", "Adm", "Operaciones", "2024-01-29", "13:6:27", "10:54:04", "13:06:27"),
    @("15961357", "Mauricio Sanchez", "Administrativa", "2024-01-29", "13:6:30", "10:54:06", "13:06:30"),
    @("1054398414", "Julian Largo", "Administrativa", "2024-01-29", "13:6:31", "10:54:08", "13:06:31")
)

$startRow = 5
$endRow = $startRow + $data.Count - 1

# Columns A (ID) and D (date) contain digit-only / date-shaped text that
# Excel would otherwise auto-convert to a number or a date serial; force
# them to stay text so the stored value matches the original literal string.
$ws.Range("A$startRow`:A$endRow").NumberFormat = "@"
$ws.Range("D$startRow`:D$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]
    for ($c = 0; $c -lt $rowData.Count; $c++) {
        $ws.Cells.Item($row, $c + 1).Value = $rowData[$c]
    }
}
